# Tag-1_1-Begruessung.pptx edit:
#  1. Update the date placeholder on the slide master from 13.06.2024 to 14.06.2024
#  2. Remove the "Rectangle 39" shape (the Brockhaus AG hyperlinked rectangle) from the slide master

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# 1) Fix the date field text on "Rectangle 6"
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 6") {
        $shp.TextFrame.TextRange.Text = "14.06.2024"
    }
}

# 2) Remove "Rectangle 39" (Brockhaus AG hyperlink) shape from the slide master
for ($i = $m.Shapes.Count; $i -ge 1; $i--) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 39") {
        $shp.Delete()
    }
}
